# Change the "Unsupported" label on slide 1 to "Refuted"
# (commit: "Change 'Supported' to 'Refuted'")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.TextRange.Text -eq "Unsupported") {
            $shape.TextFrame.TextRange.Text = "Refuted"
        }
    }
}
